$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure cells keep their text formatting (values look numeric/percent but are stored as text)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.749.53'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.744.02'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -5.02%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.93'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -8.99%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5038'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -6.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.98'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2692'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -10.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06146'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -10.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.744.70'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06923'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.45'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -12.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.522'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6006'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -18.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '76.89'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -13.51%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.749.72'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006851'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -13.32%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -16.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.965.91'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.040'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -12.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.249'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -12.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.149'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -11.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '137.59'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.514'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -11.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '14.99'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -11.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.810'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -17.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '104.02'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -6.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08117'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -8.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.746'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -11.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.471'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -14.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04573'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.619'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -10.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9831'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -13.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6106'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -16.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.668'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -13.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01552'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -9.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.916'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -14.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.67'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3823'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -18.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.065'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -14.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7331'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -18.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05370'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -6.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1111'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -10.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.942'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -19.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.19'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -13.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.45'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -12.66%  '
